# Refined metadata to be additional tab
#
# 1) Refresh the "query time" timestamps recorded in column F of the
#    existing "data" sheet (these are regenerated every time the panel
#    data is re-fetched/exported).
# 2) Add a new "metadata" worksheet (after "data") capturing the panel
#    query metadata that used to live only implicitly - now its own tab.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)
$data.Name = "data"

# --- 1) Update the per-row query timestamps on the "data" sheet -----------
$timestamps = @{
    2  = "2021-10-05 14:22:37.130472"
    3  = "2021-10-05 14:22:37.130480"
    4  = "2021-10-05 14:22:37.130483"
    5  = "2021-10-05 14:22:37.130486"
    6  = "2021-10-05 14:22:37.130488"
    7  = "2021-10-05 14:22:37.130491"
    8  = "2021-10-05 14:22:37.130493"
    9  = "2021-10-05 14:22:37.130496"
    10 = "2021-10-05 14:22:37.130499"
    11 = "2021-10-05 14:22:37.130501"
    12 = "2021-10-05 14:22:37.130503"
    13 = "2021-10-05 14:22:37.130506"
    14 = "2021-10-05 14:22:37.130508"
    15 = "2021-10-05 14:22:37.130511"
    16 = "2021-10-05 14:22:37.130513"
    17 = "2021-10-05 14:22:37.130516"
    18 = "2021-10-05 14:22:37.130518"
    19 = "2021-10-05 14:22:37.130521"
    20 = "2021-10-05 14:22:37.130523"
    21 = "2021-10-05 14:22:37.130526"
    22 = "2021-10-05 14:22:37.130528"
    23 = "2021-10-05 14:22:37.130531"
    24 = "2021-10-05 14:22:37.130533"
    25 = "2021-10-05 14:22:37.130535"
    26 = "2021-10-05 14:22:37.130538"
    27 = "2021-10-05 14:22:37.130541"
    28 = "2021-10-05 14:22:37.130543"
    29 = "2021-10-05 14:22:37.130546"
    30 = "2021-10-05 14:22:37.130548"
    31 = "2021-10-05 14:22:37.130550"
    32 = "2021-10-05 14:22:37.130553"
    33 = "2021-10-05 14:22:37.130555"
    34 = "2021-10-05 14:22:37.130558"
}

foreach ($row in $timestamps.Keys) {
    $data.Cells.Item($row, 6).Value = $timestamps[$row]
}

# --- 2) Add the new "metadata" worksheet -----------------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$headers = @{
    2 = "data_name"
    3 = "data_id"
    4 = "data_version"
    5 = "data_version_created"
    6 = "panel_query_time"
    7 = "panel_get_request"
}
foreach ($col in $headers.Keys) {
    $meta.Cells.Item(1, $col).Value = $headers[$col]
}

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Sarcoma cancer susceptibility"
$meta.Cells.Item(2, 3).Value = 217
# "1.20" is a version label, not a number - force text so the trailing
# zero (and type) survive, matching the source export. Apply the "@"
# format just long enough to enter the value, then reset the cell back
# to the sheet's plain/default formatting (copied from an unstyled data
# cell) so the text-format style doesn't linger on the saved cell.
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.20"
$data.Range("B2").Copy() | Out-Null
$meta.Cells.Item(2, 4).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$meta.Cells.Item(2, 5).Value = "2021-02-11T18:05:06.278877Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:22:37.126855"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/217/?format=json"

# Match the bold/border/centered header style used on the "data" sheet
# (copy the formatting only, so the values entered above are preserved).
$data.Range("B1").Copy() | Out-Null
$meta.Range("B1:G1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$data.Range("A2").Copy() | Out-Null
$meta.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$meta.Range("A1").Select() | Out-Null
$data.Select() | Out-Null
